$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Results" values from D2 and D4 (column D data dropped)
$ws.Range("D2").ClearContents()
$ws.Range("D4").ClearContents()

# Add a new row of test data: TC4 / test888@mail.com / test123
$ws.Range("A5").Value = "TC4"
$ws.Range("B5").Value = "test888@mail.com"
$ws.Range("C5").Value = "test123"

# Turn the new UserName cell into a mailto hyperlink like the others
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:test888@mail.com")

# Match the style used by the other UserName hyperlink cells
$ws.Range("B5").Style = $ws.Range("B4").Style

# Move the active selection as in the edited workbook
$ws.Range("D9").Select() | Out-Null
